$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37.15942529324055
$ws.Range("C2").Value = 38.474087974276
$ws.Range("D2").Value = 35.83432145104847
$ws.Range("E2").Value = 37.16020833104659
$ws.Range("F2").Value = 37.30451749110371
$ws.Range("G2").Value = 37.28936261909382
$ws.Range("H2").Value = 40.33496757003094
$ws.Range("I2").Value = 32.91175718211949
$ws.Range("J2").Value = 37.28202344561605
$ws.Range("K2").Value = 37.05413863822732
$ws.Range("L2").Value = 37.17916710262138
$ws.Range("M2").Value = 37.28251039841725
$ws.Range("N2").Value = 20.34401335102639
$ws.Range("O2").Value = 32.52730289766801
$ws.Range("P2").Value = 41.33733540220842
$ws.Range("Q2").Value = 34.40099778076166
